$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
# covering 2021-09-02 through 2021-09-09 ("aggiornamento a 9/09 compreso")
$data = @(
    @(44441, 1, 13, 153.1754447979262),
    @(44442, 0, 13, 153.1754447979262),
    @(44443, 0, 9, 106.0445387062566),
    @(44444, 0, 4, 47.13090609166961),
    @(44445, 0, 1, 11.7827265229174),
    @(44446, 1, 2, 23.5654530458348),
    @(44447, 0, 2, 23.5654530458348),
    @(44448, 3, 4, 47.13090609166961)
)

$lastRow = 366
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value2 = $vals[0]
    $ws.Cells.Item($row, 2).Value2 = $vals[1]
    $ws.Cells.Item($row, 3).Value2 = $vals[2]
    $ws.Cells.Item($row, 4).Value2 = $vals[3]

    # Column A carries the date style (bold, centered, bordered, custom
    # date/time number format) used throughout the column; replicate it
    # onto the newly appended cell without disturbing the value just set.
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
